$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.833.45'
$ws.Range('E2').Value = '  +0.41%  '

# Row 3: update D3, E3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.633.58'
$ws.Range('E3').Value = '  -0.06%  '

# Row 4: update D4, E4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.44%  '

# Row 5: update D5, E5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.40'
$ws.Range('E5').Value = '  -0.43%  '

# Row 6: update E6
$ws.Range('E6').Value = '  -0.20%  '

# Row 7: update D7, E7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.47%  '

# Row 8: update E8
$ws.Range('E8').Value = '  -0.73%  '

# Row 9: update D9, E9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0631'
$ws.Range('E9').Value = '  -0.57%  '

# Row 10: update D10, E10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.66'
$ws.Range('E10').Value = '  +0.60%  '

# Row 11: update D11, E11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0790'
$ws.Range('E11').Value = '  +0.54%  '

# Row 12: update D12, E12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.27'
$ws.Range('E12').Value = '  +0.74%  '

# Row 13: update D13, E13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.856.39'
$ws.Range('E13').Value = '  -0.21%  '

# Row 14: update D14, E14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.624.67'
$ws.Range('E14').Value = '  -0.83%  '

# Row 15: update D15, E15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.552'
$ws.Range('E15').Value = '  -0.55%  '

# Row 16: update D16, E16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₃0758'
$ws.Range('E16').Value = '  -0.74%  '

# Row 17: update D17, E17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.77'
$ws.Range('E17').Value = '  +0.10%  '

# Row 18: update D18, E18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.815.27'
$ws.Range('E18').Value = '  +0.23%  '

# Row 19: update E19
$ws.Range('E19').Value = '  -0.41%  '

# Row 20: update E20
$ws.Range('E20').Value = '  +0.04%  '

# Row 21: update D21, E21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.26'
$ws.Range('E21').Value = '  -1.19%  '

# Row 22: update D22, E22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.95'
$ws.Range('E22').Value = '  +0.25%  '

# Row 23: update D23, E23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.32'
$ws.Range('E23').Value = '  +0.89%  '

# Row 24: update D24, E24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.42%  '

# Row 25: update E25
$ws.Range('E25').Value = '  -1.86%  '

# Row 26: update D26, E26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.49'
$ws.Range('E26').Value = '  +1.45%  '

# Row 27: update D27, E27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.123'
$ws.Range('E27').Value = '  +1.06%  '

# Row 28: update E28
$ws.Range('E28').Value = '  -0.62%  '

# Row 29: update D29, E29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.49'
$ws.Range('E29').Value = '  -0.03%  '

# Row 30: update E30
$ws.Range('E30').Value = '  -0.61%  '

# Row 31: update E31
$ws.Range('E31').Value = '  +0.43%  '

# Row 32: update E32
$ws.Range('E32').Value = '  -0.25%  '

# Row 33: update D33, E33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.23'
$ws.Range('E33').Value = '  -0.38%  '

# Row 34: update E34
$ws.Range('E34').Value = '  +0.59%  '

# Row 35: update E35
$ws.Range('E35').Value = '  +0.25%  '

# Row 36: update D36, E36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.905'
$ws.Range('E36').Value = '  +0.72%  '

# Row 37: update D37, E37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.145.52'
$ws.Range('E37').Value = '  +2.05%  '

# Row 38: update D38, E38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.545'
$ws.Range('E38').Value = '  -0.28%  '

# Row 39: update D39, E39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.50'
$ws.Range('E39').Value = '  -1.17%  '

# Row 40: update E40
$ws.Range('E40').Value = '  +0.30%  '

# Row 41: update D41, E41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.997'
$ws.Range('E41').Value = '  -0.54%  '

# Row 42: update E42
$ws.Range('E42').Value = '  +0.71%  '

# Row 43: update D43, E43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.50'
$ws.Range('E43').Value = '  +0.86%  '

# Row 44: update E44
$ws.Range('E44').Value = '  +0.03%  '

# Row 45: update D45, E45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.766.28'
$ws.Range('E45').Value = '  -0.18%  '

# Row 46: update B46, C46, D46, E46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '55.50'
$ws.Range('E46').Value = '  +0.90%  '

# Row 47: update B47, C47, D47, E47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.0512'
$ws.Range('E47').Value = '  +2.27%  '

# Row 48: update B48, C48, D48, E48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '1.46'
$ws.Range('E48').Value = '  +5.68%  '

# Row 49: update B49, C49, D49, E49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '0.416'
$ws.Range('E49').Value = '  -0.18%  '

# Row 50: update B50, C50, D50, E50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.56'
$ws.Range('E50').Value = '  -0.40%  '

# Row 51: update B51, C51, D51, E51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '0.0958'
$ws.Range('E51').Value = '  +2.11%  '
